# edit.ps1 - reproduces the "revert the merge" edit on MarketGap.pptx
#
# Summary of the target change (see diff):
#  - Chart1 (bubble chart) on slide 1:
#      * remove the major gridlines from both value axes
#      * swap the axis-line formatting: the horizontal ("Build Area") axis line
#        picks up the old gridline styling (thin light-gray line), while the
#        vertical ("Cost") axis line becomes invisible (noFill)
#      * move the legend's manual layout position
#  - The cached display text of the (automatically updating) date fields on
#    every slide layout / the slide master / the notes master is restored
#    from "11/16/2014" back to "2014-11-16"
#
# Internal/non-semantic bookkeeping values that PowerPoint regenerates on its
# own (chart axis IDs, the per-slide p14:modId coauthoring checksum, and the
# physical ordering of the customXml parts inside the zip package) are not
# exposed anywhere in the PowerPoint object model, so they are intentionally
# left alone here.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Chart formatting (slide 1, shape 4 = "Chart 8")
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$chartShape = $slide.Shapes.Item(4)
$chart = $chartShape.Chart

# Value axis 1 = horizontal ("Build Area (mm2)") axis
$axBuildArea = $chart.Axes(1)
# Value axis 2 = vertical ("Cost ($)") axis
$axCost = $chart.Axes(2)

# Drop the major gridlines on both value axes.
$axBuildArea.HasMajorGridlines = $false
$axCost.HasMajorGridlines = $false

# The horizontal axis line now uses the thin light-gray styling that used to
# belong to the gridlines (0.75pt / light gray).
$buildAreaLine = $axBuildArea.Format.Line
$buildAreaLine.Weight = 0.75
$buildAreaLine.ForeColor.RGB = 14277081

# The vertical axis line becomes invisible.
$axCost.Format.Line.Visible = $false

# Reposition the legend (manual layout, expressed as a fraction of the chart
# frame -> convert the target fractions into absolute points using the
# chart's on-slide size).
$legend = $chart.Legend
$legend.Left = 478.11992125984256
$legend.Top = 30.96937007874016

# ---------------------------------------------------------------------
# 2. Restore the cached date-field text across every layout / the master /
#    the notes master (field stays an auto date field; only the last-cached
#    display string changes).
# ---------------------------------------------------------------------
function Set-DatePlaceholderText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "2014-11-16"
        }
    }
}

Set-DatePlaceholderText $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes
}

Set-DatePlaceholderText $p.NotesMaster.Shapes
